$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "`"type`": [ `"string`" ],`n`"format`": `"string`",`n`"pattern`": `"^([{]\s*([0-9]{0,})((\s*[,]\s*[0-9]{1,}){0,})\s*[}])$`",`n`"minLength`": 2"
$ws.Range("D5").Value = "{}`n{123}`n{123,456,789}`n{ 123 , 456 , 789 }"
$ws.Range("E5").Value = "2`n5`n13`n19"
